$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (WP4217 / Ebola Virus Pathway on Host)
$ws.Range("E2").Value = -2.72793281063790216123265963688
$ws.Range("F2").Value = 0.000032516350321289133606939004
$ws.Range("G2").Value = 0.000780392407710939206566536086
$ws.Range("H2").Value = 0.000753010217966695719928049702

# Row 3 (WP2328 / Allograft Rejection)
$ws.Range("E3").Value = -2.297730064065564015152176580159
$ws.Range("F3").Value = 0.001900620765738523371135570983
$ws.Range("G3").Value = 0.022807449188862280453626851795
$ws.Range("H3").Value = 0.022007187813814482052476861895
